# Auto-generated edit script applying cryptos.xlsx data refresh
# (commit: "Updated cryptos list ... with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Subscript-three character used in PEPE's price (0.0<sub>3</sub>0869)
$sub3 = [char]0x2083

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '72.900.83'
$ws.Cells.Item(2, 5).Value = '  +0.88%  '

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '4.019.98'
$ws.Cells.Item(3, 5).Value = '  -0.21%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.05%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '592.61'
$ws.Cells.Item(5, 5).Value = '  +12.01%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '153.04'
$ws.Cells.Item(6, 5).Value = '  +1.22%  '

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '0.689'
$ws.Cells.Item(7, 5).Value = '  -1.82%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  -0.08%  '

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.761'
$ws.Cells.Item(9, 5).Value = '  +1.45%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '0.171'
$ws.Cells.Item(10, 5).Value = '  -0.75%  '

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '54.56'
$ws.Cells.Item(11, 5).Value = '  +9.55%  '

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '11.00'
$ws.Cells.Item(13, 5).Value = '  +3.04%  '

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '4.659.58'
$ws.Cells.Item(14, 5).Value = '  -0.20%  '

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '4.034.57'
$ws.Cells.Item(15, 5).Value = '  +0.26%  '

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '1.28'
$ws.Cells.Item(16, 5).Value = '  +7.58%  '

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '14.28'
$ws.Cells.Item(17, 5).Value = '  +0.93%  '

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '20.71'
$ws.Cells.Item(18, 5).Value = '  +0.49%  '

# Row 19
$ws.Cells.Item(19, 5).Value = '  -0.48%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '72.750.87'
$ws.Cells.Item(20, 5).Value = '  +0.92%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '439.14'
$ws.Cells.Item(21, 5).Value = '  +1.50%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '4.81'
$ws.Cells.Item(22, 5).Value = '  +13.12%  '

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '97.26'
$ws.Cells.Item(23, 5).Value = '  -0.78%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '3.53'
$ws.Cells.Item(24, 5).Value = '  +0.73%  '

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '14.42'
$ws.Cells.Item(25, 5).Value = '  +0.97%  '

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '4.36'
$ws.Cells.Item(26, 5).Value = '  +19.91%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '11.45'
$ws.Cells.Item(27, 5).Value = '  -0.19%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '10.79'
$ws.Cells.Item(28, 5).Value = '  +0.36%  '

# Row 29
$ws.Cells.Item(29, 5).Value = '  +1.22%  '

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '36.75'
$ws.Cells.Item(30, 5).Value = '  -0.19%  '

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '7.99'
$ws.Cells.Item(31, 5).Value = '  +9.61%  '

# Row 32
$ws.Cells.Item(32, 2).Value = 'Hedera'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '0.134'
$ws.Cells.Item(32, 5).Value = '  +1.79%  '

# Row 33
$ws.Cells.Item(33, 2).Value = 'Cosmos'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '13.68'
$ws.Cells.Item(33, 5).Value = '  +1.33%  '

# Row 34
$ws.Cells.Item(34, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '50.37'
$ws.Cells.Item(34, 5).Value = '  +5.27%  '

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '686.09'
$ws.Cells.Item(35, 5).Value = '  +0.72%  '

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '71.38'
$ws.Cells.Item(36, 5).Value = '  +8.18%  '

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '0.445'
$ws.Cells.Item(37, 5).Value = '  -1.44%  '

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = "0.0${sub3}0869"
$ws.Cells.Item(38, 5).Value = '  +4.35%  '

# Row 39
$ws.Cells.Item(39, 2).Value = 'Kaspa'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '0.149'
$ws.Cells.Item(39, 5).Value = '  -1.27%  '

# Row 40
$ws.Cells.Item(40, 2).Value = 'WEMIXToken'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '3.41'
$ws.Cells.Item(40, 5).Value = '  +3.92%  '

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '11.22'
$ws.Cells.Item(41, 5).Value = '  +10.50%  '

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '3.35'
$ws.Cells.Item(42, 5).Value = '  -1.96%  '

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '1.00'

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '0.0494'
$ws.Cells.Item(44, 5).Value = '  +0.22%  '

# Row 45
$ws.Cells.Item(45, 5).Value = '  +0.12%  '

# Row 46
$ws.Cells.Item(46, 2).Value = 'Fetch.AI'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '2.78'
$ws.Cells.Item(46, 5).Value = '  +1.44%  '

# Row 47
$ws.Cells.Item(47, 2).Value = 'Stellar'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '0.151'
$ws.Cells.Item(47, 5).Value = '  +0.23%  '

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '3.41'
$ws.Cells.Item(48, 5).Value = '  +1.05%  '

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '3.52'
$ws.Cells.Item(49, 5).Value = '  +7.71%  '

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '3.04'
$ws.Cells.Item(50, 5).Value = '  -0.50%  '

# Row 51
$ws.Cells.Item(51, 2).Value = 'Maker'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '2.839.81'
$ws.Cells.Item(51, 5).Value = '  +11.24%  '
